$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.662.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.526.63"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.14"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.914.49"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.534.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.20"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.809"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.647.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.57"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("E24").Value = "  -1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -3.48%  "

$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.53"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("E33").Value = "  -2.26%  "

$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.65%  "

$ws.Range("E37").Value = "  -4.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.85%  "

$ws.Range("E39").Value = "  -2.99%  "

$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.80"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.027.83"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0295"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("E46").Value = "  -4.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.83"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.770.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.72%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.187"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.81%  "
